$wb = $excel.ActiveWorkbook

# Update values on the PSA_LOLO sheet
$psa = $wb.Worksheets.Item("PSA_LOLO")
$psa.Range("B2").Value = 32214
$psa.Range("B3").Value = 8142

# Make PSA_LOLO the active/selected sheet (tabSelected moves from
# OverallRebateEfficiency to PSA_LOLO)
$psa.Activate()

$wb.Save()
